$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 116.5
$ws.Cells.Item(5, 9).Value = 133.33333
$ws.Cells.Item(5, 11).Value = 133.33333
$ws.Cells.Item(5, 13).Value = -18.33332999999999

$ws.Cells.Item(18, 8).Value = 400.2857
$ws.Cells.Item(18, 9).Value = 400.2857
$ws.Cells.Item(18, 11).Value = 400.2857
$ws.Cells.Item(18, 13).Value = -116.2857

$ws.Cells.Item(28, 8).Value = 1252.5
$ws.Cells.Item(28, 9).Value = 1252.5
$ws.Cells.Item(28, 11).Value = 1252.5
$ws.Cells.Item(28, 13).Value = -767.5

$ws.Cells.Item(32, 8).Value = 10345.5
$ws.Cells.Item(32, 10).Value = 10062.125
$ws.Cells.Item(32, 12).Value = 10062.125
$ws.Cells.Item(32, 14).Value = -10714.125

$ws.Cells.Item(62, 8).Value = 3310.889
$ws.Cells.Item(62, 9).Value = 2849.75
$ws.Cells.Item(62, 11).Value = 2849.75
$ws.Cells.Item(62, 13).Value = -2225.75

$ws.Cells.Item(65, 8).Value = 3310.889
$ws.Cells.Item(65, 9).Value = 2849.75
$ws.Cells.Item(65, 11).Value = 14248.75
$ws.Cells.Item(65, 13).Value = -11128.75

$ws.Cells.Item(106, 8).Value = 3454.5
$ws.Cells.Item(106, 9).Value = 3454.5
$ws.Cells.Item(106, 11).Value = 3454.5
$ws.Cells.Item(106, 13).Value = -2823.5

$ws.Cells.Item(129, 8).Value = 1979.949
$ws.Cells.Item(129, 10).Value = 2096.8242
$ws.Cells.Item(129, 12).Value = 6290.4726
$ws.Cells.Item(129, 14).Value = -16290.4726

$ws.Cells.Item(138, 8).Value = 3243.5715
$ws.Cells.Item(138, 9).Value = 1499.5
$ws.Cells.Item(138, 10).Value = 3427.158
$ws.Cells.Item(138, 11).Value = 4498.5
$ws.Cells.Item(138, 12).Value = 10281.474
$ws.Cells.Item(138, 13).Value = 641.5
$ws.Cells.Item(138, 14).Value = -20561.474

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3385.5
$ws.Cells.Item(2, 9).Value = 3183.0833
$ws.Cells.Item(2, 11).Value = 3183.0833
$ws.Cells.Item(2, 13).Value = -3070.0833

$ws.Cells.Item(4, 8).Value = 185.6
$ws.Cells.Item(4, 9).Value = 195.75
$ws.Cells.Item(4, 11).Value = 195.75
$ws.Cells.Item(4, 13).Value = -79.75

$ws.Cells.Item(44, 8).Value = 42000
$ws.Cells.Item(44, 10).Value = 42000
$ws.Cells.Item(44, 12).Value = 42000
$ws.Cells.Item(44, 14).Value = -42976

$ws.Cells.Item(52, 8).Value = 100000
$ws.Cells.Item(52, 10).Value = 100000
$ws.Cells.Item(52, 12).Value = 100000
$ws.Cells.Item(52, 14).Value = -100636

$ws.Cells.Item(55, 8).Value = 23999.666
$ws.Cells.Item(55, 10).Value = 42000
$ws.Cells.Item(55, 12).Value = 42000
$ws.Cells.Item(55, 14).Value = -42630

$ws.Cells.Item(74, 8).Value = 2583.3684
$ws.Cells.Item(74, 9).Value = 2165.7778
$ws.Cells.Item(74, 11).Value = 2165.7778
$ws.Cells.Item(74, 13).Value = -1291.7778

$ws.Cells.Item(77, 8).Value = 2583.3684
$ws.Cells.Item(77, 9).Value = 2165.7778
$ws.Cells.Item(77, 11).Value = 10828.889
$ws.Cells.Item(77, 13).Value = -6460.888999999999

$ws.Cells.Item(116, 8).Value = 3385.5
$ws.Cells.Item(116, 9).Value = 3183.0833
$ws.Cells.Item(116, 11).Value = 3183.0833
$ws.Cells.Item(116, 13).Value = -889.0832999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3385.5
$ws.Cells.Item(3, 9).Value = 3183.0833
$ws.Cells.Item(3, 11).Value = 3183.0833
$ws.Cells.Item(3, 13).Value = -3069.0833

$ws.Cells.Item(7, 8).Value = 4000
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 4000
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 4000
$ws.Cells.Item(7, 13).ClearContents()
$ws.Cells.Item(7, 14).Value = -4226

$ws.Cells.Item(20, 8).Value = 9845.909
$ws.Cells.Item(20, 9).Value = 13342.286
$ws.Cells.Item(20, 11).Value = 13342.286
$ws.Cells.Item(20, 13).Value = -13095.286

$ws.Cells.Item(102, 8).Value = 2568.6667
$ws.Cells.Item(102, 9).Value = 2568.6667
$ws.Cells.Item(102, 11).Value = 2568.6667
$ws.Cells.Item(102, 13).Value = 676.3332999999998

$ws.Cells.Item(134, 8).Value = 2642.0205
$ws.Cells.Item(134, 9).Value = 2567.319
$ws.Cells.Item(134, 11).Value = 7701.957
$ws.Cells.Item(134, 13).Value = -5166.957

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 159.35
$ws.Cells.Item(7, 9).Value = 176.47058
$ws.Cells.Item(7, 11).Value = 176.47058
$ws.Cells.Item(7, 13).Value = -63.47058000000001

$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 14).ClearContents()

$ws.Cells.Item(59, 8).Value = 23522.4
$ws.Cells.Item(59, 9).Value = 19104
$ws.Cells.Item(59, 10).Value = 24627
$ws.Cells.Item(59, 11).Value = 19104
$ws.Cells.Item(59, 12).Value = 24627
$ws.Cells.Item(59, 13).Value = -17959
$ws.Cells.Item(59, 14).Value = -26917

$ws.Cells.Item(107, 8).Value = 2090.0625
$ws.Cells.Item(107, 9).Value = 2182.7334
$ws.Cells.Item(107, 11).Value = 2182.7334
$ws.Cells.Item(107, 13).Value = -262.7334000000001

$ws.Cells.Item(111, 8).Value = 54998.832
$ws.Cells.Item(111, 10).Value = 54998.832
$ws.Cells.Item(111, 12).Value = 54998.832
$ws.Cells.Item(111, 14).Value = -63178.832

$ws.Cells.Item(129, 8).Value = 41250
$ws.Cells.Item(129, 10).Value = 41250
$ws.Cells.Item(129, 12).Value = 41250
$ws.Cells.Item(129, 14).Value = -51250

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(25, 8).Value = 3349.75
$ws.Cells.Item(25, 9).Value = 699.5
$ws.Cells.Item(25, 10).Value = 6000
$ws.Cells.Item(25, 11).Value = 2098.5
$ws.Cells.Item(25, 12).Value = 18000
$ws.Cells.Item(25, 13).Value = -1929.5
$ws.Cells.Item(25, 14).Value = -18338

$ws.Cells.Item(30, 8).Value = 3349.75
$ws.Cells.Item(30, 9).Value = 699.5
$ws.Cells.Item(30, 10).Value = 6000
$ws.Cells.Item(30, 11).Value = 2098.5
$ws.Cells.Item(30, 12).Value = 18000
$ws.Cells.Item(30, 13).Value = -1996.5
$ws.Cells.Item(30, 14).Value = -18204

$ws.Cells.Item(99, 8).Value = 5666.5
$ws.Cells.Item(99, 9).Value = 1333
$ws.Cells.Item(99, 11).Value = 3999
$ws.Cells.Item(99, 13).Value = -1753

$ws.Cells.Item(115, 8).Value = 0
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 14).ClearContents()

$ws.Cells.Item(138, 8).Value = 4285.5
$ws.Cells.Item(138, 9).Value = 4285.5
$ws.Cells.Item(138, 11).Value = 12856.5
$ws.Cells.Item(138, 13).Value = -7716.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 83.56521600000001
$ws.Cells.Item(2, 9).Value = 67.7619
$ws.Cells.Item(2, 11).Value = 67.7619
$ws.Cells.Item(2, 13).Value = 45.2381

$ws.Cells.Item(126, 8).Value = 8186.4
$ws.Cells.Item(126, 9).Value = 3571
$ws.Cells.Item(126, 10).Value = 12224.875
$ws.Cells.Item(126, 11).Value = 10713
$ws.Cells.Item(126, 12).Value = 36674.625
$ws.Cells.Item(126, 13).Value = -8243
$ws.Cells.Item(126, 14).Value = -41614.625

$ws.Cells.Item(132, 8).Value = 3000.4
$ws.Cells.Item(132, 9).Value = 2823.3635
$ws.Cells.Item(132, 11).Value = 8470.0905
$ws.Cells.Item(132, 13).Value = -5940.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 850.2778
$ws.Cells.Item(16, 10).Value = 219
$ws.Cells.Item(16, 12).Value = 219
$ws.Cells.Item(16, 14).Value = -559

$ws.Cells.Item(24, 8).Value = 9000
$ws.Cells.Item(24, 10).Value = 9000
$ws.Cells.Item(24, 12).Value = 9000
$ws.Cells.Item(24, 14).Value = -9686

$ws.Cells.Item(43, 8).Value = 32499.75
$ws.Cells.Item(43, 9).Value = 78999
$ws.Cells.Item(43, 10).Value = 17000
$ws.Cells.Item(43, 11).Value = 78999
$ws.Cells.Item(43, 12).Value = 17000
$ws.Cells.Item(43, 13).Value = -78806
$ws.Cells.Item(43, 14).Value = -17386

$ws.Cells.Item(55, 8).Value = 133.66667
$ws.Cells.Item(55, 10).Value = 174.5
$ws.Cells.Item(55, 12).Value = 174.5
$ws.Cells.Item(55, 14).Value = -520.5

$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 14).ClearContents()

$ws.Cells.Item(129, 8).Value = 80000
$ws.Cells.Item(129, 10).Value = 80000
$ws.Cells.Item(129, 12).Value = 80000
$ws.Cells.Item(129, 14).Value = -90000

$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 12).Value = 0
$ws.Cells.Item(130, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 5195.077
$ws.Cells.Item(132, 9).Value = 3783.1428
$ws.Cells.Item(132, 11).Value = 11349.4284
$ws.Cells.Item(132, 13).Value = -8819.428400000001

$ws.Cells.Item(136, 8).Value = 83347720
$ws.Cells.Item(136, 9).Value = 15965.8
$ws.Cells.Item(136, 11).Value = 47897.39999999999
$ws.Cells.Item(136, 13).Value = -45347.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2556.0625
$ws.Cells.Item(81, 10).Value = 2785.5715
$ws.Cells.Item(81, 12).Value = 5571.143
$ws.Cells.Item(81, 14).Value = -7693.143

$ws.Cells.Item(84, 8).Value = 2556.0625
$ws.Cells.Item(84, 10).Value = 2785.5715
$ws.Cells.Item(84, 12).Value = 27855.715
$ws.Cells.Item(84, 14).Value = -38463.715

$ws.Cells.Item(86, 8).Value = 40000
$ws.Cells.Item(86, 10).Value = 40000
$ws.Cells.Item(86, 12).Value = 40000
$ws.Cells.Item(86, 14).Value = -42246

$ws.Cells.Item(89, 8).Value = 40000
$ws.Cells.Item(89, 10).Value = 40000
$ws.Cells.Item(89, 12).Value = 200000
$ws.Cells.Item(89, 14).Value = -211232

$ws.Cells.Item(132, 8).Value = 3119.4211
$ws.Cells.Item(132, 9).Value = 1822.4348
$ws.Cells.Item(132, 11).Value = 5467.3044
$ws.Cells.Item(132, 13).Value = -2937.3044
